$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 03:03:19"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 03:03:14"
$wsZhCn.Range("K2").Value = "2016-08-27 03:03:42"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-27 03:03:19"
$wsDeDe.Range("K2").Value = "2016-08-27 03:03:49"
